$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.Value = "'" + $text
}

$ws.Range("D2").Value = "61.289.41"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "2.923.72"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "597.63"

Set-TextValue $ws.Range("D6") "144.74"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.24%  "

Set-TextValue $ws.Range("D9") "6.93"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("E12").Value = "  -1.10%  "

Set-TextValue $ws.Range("D13") "33.43"
$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "3.407.73"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "61.250.97"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "2.925.69"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("E18").Value = "  -0.63%  "

Set-TextValue $ws.Range("D19") "431.79"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("E21").Value = "  -1.31%  "

Set-TextValue $ws.Range("D22") "7.07"
$ws.Range("E22").Value = "  -0.13%  "

Set-TextValue $ws.Range("D23") "81.79"
$ws.Range("E23").Value = "  +0.27%  "

Set-TextValue $ws.Range("D24") "10.90"
$ws.Range("E24").Value = "  -1.12%  "

Set-TextValue $ws.Range("D25") "2.17"
$ws.Range("E25").Value = "  -2.07%  "

Set-TextValue $ws.Range("D26") "11.74"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("E27").Value = "  +0.01%  "

Set-TextValue $ws.Range("D28") "2.19"
$ws.Range("E28").Value = "  -5.08%  "

$ws.Range("E29").Value = "  -0.86%  "

Set-TextValue $ws.Range("D30") "6.88"
$ws.Range("E30").Value = "  -2.74%  "

Set-TextValue $ws.Range("D31") "26.61"
$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "0.0₃0879"
$ws.Range("E34").Value = "  +3.20%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("E39").Value = "  -1.07%  "

Set-TextValue $ws.Range("D40") "8.55"
$ws.Range("E40").Value = "  -0.29%  "

Set-TextValue $ws.Range("D41") "42.52"
$ws.Range("E41").Value = "  +5.09%  "

Set-TextValue $ws.Range("D42") "0.279"
$ws.Range("E42").Value = "  -2.71%  "

$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "2.693.54"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D45") "133.73"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D46") "365.27"
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("E50").Value = "  -1.21%  "

Set-TextValue $ws.Range("D51") "0.125"
$ws.Range("E51").Value = "  -0.84%  "
